$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7801687900233389
$ws.Range("C2").Value = 0.04160859039693321
$ws.Range("D2").Value = 0.143342175087998
$ws.Range("E2").Value = 0.06113185119995279
$ws.Range("F2").Value = 2.630708748447447
$ws.Range("I2").Value = 2.147253301571425
$ws.Range("K2").Value = 0.6164355302459512
$ws.Range("L2").Value = 0.2456817495575621
$ws.Range("B3").Value = 0.7541894255479917
$ws.Range("C3").Value = 0.03620255580473497
$ws.Range("D3").Value = 0.1426570395144076
$ws.Range("E3").Value = 0.06068365956370769
$ws.Range("F3").Value = 2.577096267988736
$ws.Range("I3").Value = 2.114043808880069
$ws.Range("K3").Value = 0.5853739611387425
$ws.Range("L3").Value = 0.238526883288003
$ws.Range("B4").Value = 0.7388538796103887
$ws.Range("C4").Value = 0.03288011789013012
$ws.Range("D4").Value = 0.1422241723569364
$ws.Range("E4").Value = 0.06043441177349784
$ws.Range("F4").Value = 2.545016477067776
$ws.Range("I4").Value = 2.09420565580325
$ws.Range("K4").Value = 0.5667910953254705
$ws.Range("L4").Value = 0.2342796997543815
$ws.Range("B5").Value = 0.7327593670990211
$ws.Range("C5").Value = 0.03152531484467147
$ws.Range("D5").Value = 0.1420446831958699
$ws.Range("E5").Value = 0.06033936952587915
$ws.Range("F5").Value = 2.532153828972611
$ws.Range("I5").Value = 2.086259911447826
$ws.Range("K5").Value = 0.5593412042661896
$ws.Range("L5").Value = 0.2325855804221817
$ws.Range("B6").Value = 0.7317567298549648
$ws.Range("C6").Value = 0.03130029460496075
$ws.Range("D6").Value = 0.1420146916058158
$ws.Range("E6").Value = 0.06032398226020952
$ws.Range("F6").Value = 2.530030666377584
$ws.Range("I6").Value = 2.084948875163789
$ws.Range("K6").Value = 0.5581115653201039
$ws.Range("L6").Value = 0.2323064847839476
$ws.Range("B7").Value = 0.7387710599403192
$ws.Range("C7").Value = 0.03286185024772692
$ws.Range("D7").Value = 0.1422217642595385
$ws.Range("E7").Value = 0.06043310356150933
$ws.Range("F7").Value = 2.544842157063854
$ws.Range("I7").Value = 2.094097936746337
$ws.Range("K7").Value = 0.5666901266275488
$ws.Range("L7").Value = 0.2342567039733723
$ws.Range("B8").Value = 0.7710832498189859
$ws.Range("C8").Value = 0.03974517324338933
$ws.Range("D8").Value = 0.1431084544871055
$ws.Range("E8").Value = 0.06097193249065036
$ws.Range("F8").Value = 2.612048684312384
$ws.Range("I8").Value = 2.135687664364212
$ws.Range("K8").Value = 0.6056238609949673
$ws.Range("L8").Value = 0.2431844276398181
$ws.Range("B9").Value = 0.8393404167546521
$ws.Range("C9").Value = 0.05322382392630232
$ws.Range("D9").Value = 0.1447517514737484
$ws.Range("E9").Value = 0.06223435567556734
$ws.Range("F9").Value = 2.750534342411328
$ws.Range("I9").Value = 2.221657233728934
$ws.Range("K9").Value = 0.6858677430424507
$ws.Range("L9").Value = 0.2618532966541807
$ws.Range("B10").Value = 0.8924880890712359
$ws.Range("C10").Value = 0.06312273449840688
$ws.Range("D10").Value = 0.1459026043698088
$ws.Range("E10").Value = 0.0632874026881467
$ws.Range("F10").Value = 2.856429242130105
$ws.Range("I10").Value = 2.287554572073958
$ws.Range("K10").Value = 0.7472251482336389
$ws.Range("L10").Value = 0.2762844874469295
$ws.Range("B11").Value = 0.9173216697800797
$ws.Range("C11").Value = 0.06762696470234175
$ws.Range("D11").Value = 0.1464142668620028
$ws.Range("E11").Value = 0.06379375755389205
$ws.Range("F11").Value = 2.905520194086591
$ws.Range("I11").Value = 2.318137281930518
$ws.Range("K11").Value = 0.7756664249610026
$ws.Range("L11").Value = 0.2830065516543954
$ws.Range("B12").Value = 0.9268200931948911
$ws.Range("C12").Value = 0.06933289356190642
$ws.Range("D12").Value = 0.1466063422340866
$ws.Range("E12").Value = 0.06398942885917691
$ws.Range("F12").Value = 2.924242794660529
$ws.Range("I12").Value = 2.329805929395761
$ws.Range("K12").Value = 0.7865129249303209
$ws.Range("L12").Value = 0.2855747316200734
$ws.Range("B13").Value = 0.9247702339626755
$ws.Range("C13").Value = 0.06896547733192904
$ws.Range("D13").Value = 0.1465650496721409
$ws.Range("E13").Value = 0.063947112986245
$ws.Range("F13").Value = 2.920204623134623
$ws.Range("I13").Value = 2.327288970937929
$ws.Range("K13").Value = 0.7841735345879499
$ws.Range("L13").Value = 0.285020618295718
$ws.Range("B14").Value = 0.9181012166126834
$ws.Range("C14").Value = 0.06776730637633932
$ws.Range("D14").Value = 0.1464301025413803
$ws.Range("E14").Value = 0.06380977688741396
$ws.Range("F14").Value = 2.907057845161376
$ws.Range("I14").Value = 2.319095509045127
$ws.Range("K14").Value = 0.7765572401980876
$ws.Range("L14").Value = 0.283217382421924
$ws.Range("B15").Value = 0.9140285585466188
$ws.Range("C15").Value = 0.06703343102620352
$ws.Range("D15").Value = 0.1463472255016889
$ws.Range("E15").Value = 0.06372616572863166
$ws.Range("F15").Value = 2.899022398192443
$ws.Range("I15").Value = 2.314088206396846
$ws.Range("K15").Value = 0.7719019983118471
$ws.Range("L15").Value = 0.2821158040758007
$ws.Range("B16").Value = 0.8908783728921321
$ws.Range("C16").Value = 0.06282840504773901
$ws.Range("D16").Value = 0.1458689294891293
$ws.Range("E16").Value = 0.06325486069335895
$ws.Range("F16").Value = 2.853239590081415
$ws.Range("I16").Value = 2.285568154236543
$ws.Range("K16").Value = 0.7453771210733464
$ws.Range("L16").Value = 0.2758483547153361
$ws.Range("B17").Value = 0.8768446746979919
$ws.Range("C17").Value = 0.06024913228057471
$ws.Range("D17").Value = 0.1455724898622641
$ws.Range("E17").Value = 0.06297272529241837
$ws.Range("F17").Value = 2.825389219166084
$ws.Range("I17").Value = 2.268227540585571
$ws.Range("K17").Value = 0.7292407685753233
$ws.Range("L17").Value = 0.2720437996606506
$ws.Range("B18").Value = 0.8688346437551502
$ws.Range("C18").Value = 0.05876570611366105
$ws.Range("D18").Value = 0.1454008671597791
$ws.Range("E18").Value = 0.0628130199666046
$ws.Range("F18").Value = 2.809456859245728
$ws.Range("I18").Value = 2.258310655122528
$ws.Range("K18").Value = 0.7200093865931194
$ws.Range("L18").Value = 0.2698703164734297
$ws.Range("B19").Value = 0.8661331932640337
$ws.Range("C19").Value = 0.05826345800288379
$ws.Range("D19").Value = 0.1453425657619221
$ws.Range("E19").Value = 0.06275938822792781
$ws.Range("F19").Value = 2.804077257378367
$ws.Range("I19").Value = 2.254962737428713
$ws.Range("K19").Value = 0.716892348137435
$ws.Range("L19").Value = 0.2691369516146551
$ws.Range("B20").Value = 0.8783321911039934
$ws.Range("C20").Value = 0.06052368880453685
$ws.Range("D20").Value = 0.1456041619331216
$ws.Range("E20").Value = 0.06300249296731764
$ws.Range("F20").Value = 2.828344987096386
$ws.Range("I20").Value = 2.270067577721719
$ws.Range("K20").Value = 0.730953352645173
$ws.Range("L20").Value = 0.2724472691908915
$ws.Range("B21").Value = 0.9200575035960412
$ws.Range("C21").Value = 0.06811922991887798
$ws.Range("D21").Value = 0.1464697851736574
$ws.Range("E21").Value = 0.0638500093064458
$ws.Range("F21").Value = 2.910915758684979
$ws.Range("I21").Value = 2.321499744576826
$ws.Range("K21").Value = 0.7787922546840775
$ws.Range("L21").Value = 0.2837464204200586
$ws.Range("B22").Value = 0.9478781889630454
$ws.Range("C22").Value = 0.07308500623219061
$ws.Range("D22").Value = 0.1470257391167209
$ws.Range("E22").Value = 0.0644267911483638
$ws.Range("F22").Value = 2.965655678685948
$ws.Range("I22").Value = 2.355624719915227
$ws.Range("K22").Value = 0.8105032238776175
$ws.Range("L22").Value = 0.2912632992870101
$ws.Range("B23").Value = 0.9329793367512309
$ws.Range("C23").Value = 0.07043449129093915
$ws.Range("D23").Value = 0.146729902350728
$ws.Range("E23").Value = 0.0641168590492498
$ws.Range("F23").Value = 2.936368762678626
$ws.Range("I23").Value = 2.337364634249809
$ws.Range("K23").Value = 0.7935376308481636
$ws.Range("L23").Value = 0.2872392773052894
$ws.Range("B24").Value = 0.8776595038803805
$ws.Range("C24").Value = 0.0603995635775334
$ws.Range("D24").Value = 0.1455898467161347
$ws.Range("E24").Value = 0.06298902722036459
$ws.Range("F24").Value = 2.827008437802903
$ws.Range("I24").Value = 2.269235533546919
$ws.Range("K24").Value = 0.7301789513033725
$ws.Range("L24").Value = 0.272264817619245
$ws.Range("B25").Value = 0.8203495690275702
$ws.Range("C25").Value = 0.0495787533871237
$ws.Range("D25").Value = 0.1443172394101637
$ws.Range("E25").Value = 0.06187079730114675
$ws.Range("F25").Value = 2.712346594024069
$ws.Range("I25").Value = 2.197923114728667
$ws.Range("K25").Value = 0.6637398292864134
$ws.Range("L25").Value = 0.2566777937795166
